$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new game log numbers to the running space-separated lists ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 4 -1 1 5 6 14 3 4 4 3 2 -1 2 3 2 5 4 3 2 2"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " 14 19 10 11 7 9 1 3 8 9 5 75 6 6 25 1 5 6 2 4"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 1 2 5 -1 2 2 7 4 2 6 4 5 4 5 1 5 9 5 8 13 5 -2 4 3 0 9 1 1"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 15 16 12 18 3 12 -1 24 6 12 14"

# --- OFF sheet: update aggregate stat totals ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 441
$wsOFF.Range("F2").Value = 128
$wsOFF.Range("G2").Value = 154
$wsOFF.Range("J2").Value = 56
$wsOFF.Range("L2").Value = 661
$wsOFF.Range("M2").Value = 463
$wsOFF.Range("Q2").Value = 1213
$wsOFF.Range("B3").Value = 26
$wsOFF.Range("C3").Value = 372
$wsOFF.Range("E3").Value = 48
$wsOFF.Range("F3").Value = 230
$wsOFF.Range("H3").Value = 47
$wsOFF.Range("I3").Value = 127
$wsOFF.Range("J3").Value = 130
$wsOFF.Range("N3").Value = 38

# --- DEF sheet: update aggregate stat totals ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 388
$wsDEF.Range("E2").Value = 8
$wsDEF.Range("F2").Value = 115
$wsDEF.Range("G2").Value = 118
$wsDEF.Range("I2").Value = 9
$wsDEF.Range("J2").Value = 75
$wsDEF.Range("L2").Value = 638
$wsDEF.Range("M2").Value = 379
$wsDEF.Range("O2").Value = 58
$wsDEF.Range("Q2").Value = 1130
$wsDEF.Range("B3").Value = 20
$wsDEF.Range("C3").Value = 368
$wsDEF.Range("E3").Value = 55
$wsDEF.Range("F3").Value = 247
$wsDEF.Range("G3").Value = 77
$wsDEF.Range("H3").Value = 50
$wsDEF.Range("I3").Value = 121
$wsDEF.Range("J3").Value = 102
$wsDEF.Range("N3").Value = 57

# --- ST sheet: append to the per-kick-return / punt-return number lists ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value = $wsST.Range("B4").Value() + " 57 60 51"
$wsST.Range("B5").Value = $wsST.Range("B5").Value() + " 32 45 14"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 15 26 25"
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 42 59 36 41 57"
$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 0 6 0 7 9"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 0 0 6 5 0"

# --- ST sheet: update aggregate stat totals ---
$wsST.Range("B2").Value = 201
$wsST.Range("D2").Value = 106
$wsST.Range("F2").Value = 716
$wsST.Range("G2").Value = 699
$wsST.Range("J2").Value = 260
$wsST.Range("K2").Value = 236

# --- TURNS sheet: update aggregate stat totals ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("C2").Value = 24
$wsTURNS.Range("D2").Value = 16

# --- PEN sheet: update aggregate stat totals ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 19
